$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.679.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.368.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.87%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.360.35"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.37%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.87%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.632"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.98"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.907.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.33"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.69%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.359.91"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.86"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "64.663.35"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.987"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.09%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "460.67"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +12.28%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +9.13%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.25"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.85%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.45%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.61%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.79"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.09"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.60%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.20%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "582.31"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.06%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.34"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.141"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.72%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.93"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.44%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0759"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.93%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.106.23"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.47%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.74%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.19%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.59"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.76"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.01%  "

